$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AE1, "DiffA")
# onto the three new header cells so they match the rest of the header row.
$ws.Range("AE1").Copy()
$ws.Range("AF1:AH1").PasteSpecial(-4122)  # xlPasteFormats

# Add the three new %Diff header labels.
$ws.Range("AF1").Value = "%DiffH"
$ws.Range("AG1").Value = "%DiffD"
$ws.Range("AH1").Value = "%DiffA"
